$d = $word.ActiveDocument

function Replace-Text($findText, $replaceText) {
    $r = $d.Content
    $r.Find.Execute($findText, $true, $false, $false, $false, $false, $true, 1, $false, $replaceText, 2)
    if (-not $r.Find.Found) {
        Write-Host "NOT FOUND: $findText"
    }
}

# 1) Merge split runs back into single runs (text content is unchanged, only the
#    run-splitting is collapsed) for four bullet points.
Replace-Text "Kan søke opp band som har spilt på en scene tidligere og se nøkkelinformasjon om bandet og om tidligere konserter." "Kan søke opp band som har spilt på en scene tidligere og se nøkkelinformasjon om bandet og om tidligere konserter."

Replace-Text "Få oversikt over tidligere konserter innen en gitt sjanger, med informasjon om publikumsantall og scene." "Få oversikt over tidligere konserter innen en gitt sjanger, med informasjon om publikumsantall og scene."

Replace-Text "Kan se en rapport om konserter som viser publikumstall, kostnader og økonomisk resultat for alle konserter på en scene." "Kan se en rapport om konserter som viser publikumstall, kostnader og økonomisk resultat for alle konserter på en scene."

Replace-Text "Kan få generert et forslag til billettpris som tar høyde for markedsinformasjon og faktiske kostnader, og få forslag til billettpris på scener med ulik størrelse slik at konserter går i økonomisk balanse." "Kan få generert et forslag til billettpris som tar høyde for markedsinformasjon og faktiske kostnader, og få forslag til billettpris på scener med ulik størrelse slik at konserter går i økonomisk balanse."

# 2) Replace the "Trenger man godkjenning..." bullet with new wording, and move the
#    "_GoBack" bookmark (which Word drops at the most-recently-edited spot) here.
$target = $d.Content
$target.Find.Execute("Trenger man godkjenning for registreringen og mulighet til å slette folk fra databasen?", $true, $false, $false, $false, $false, $true, 1, $false, "Admin godkjenner nye brukere fra backend.X", 2)
if (-not $target.Find.Found) {
    Write-Host "NOT FOUND: Trenger man godkjenning..."
}

# Put the bookmark right before the temporary "X" marker (a real character, so the
# insertion point is unambiguous), then delete the marker -- this leaves the
# bookmark collapsed right after the final period, before the paragraph mark.
$markRange = $d.Content
$markRange.Find.Execute("X", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
if ($markRange.Find.Found) {
    $markRange.Collapse(1)
    if ($d.Bookmarks.Exists("_GoBack")) {
        $d.Bookmarks("_GoBack").Delete()
    }
    $d.Bookmarks.Add("_GoBack", $markRange)

    $xRange = $d.Content
    $xRange.Find.Execute("X", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
    if ($xRange.Find.Found) {
        $xRange.Text = ""
    }
}

# 3) Reword "Kan man ha flere roller?" -> "Man kan ha flere roller" (also merges the
#    split runs into one).
Replace-Text "Kan man ha flere roller?" "Man kan ha flere roller"

# 4) The "lastRenderedPageBreak" marker that used to sit on the "Booking sjef side:"
#    run moves off of it (re-writing the run via Find/Replace with identical text
#    drops the stale rendering marker).
Replace-Text "Booking sjef side:" "Booking sjef side:"
